# Update loading_percent values on Sheet1 (res_line / Case_2_230, 380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 17.76015136190095
$ws.Range("C2").Value = 10.62655083921642
$ws.Range("E2").Value = 15.75911348061313
$ws.Range("F2").Value = 37.86406441250496
$ws.Range("G2").Value = 3.647195106548639
$ws.Range("J2").Value = 7.793573481488314
$ws.Range("L2").Value = 12.66801034171153
$ws.Range("N2").Value = 17.95183981319581
$ws.Range("O2").Value = 24.21493026832341
# Row 3
$ws.Range("B3").Value = 17.31743106386717
$ws.Range("C3").Value = 10.50202157589124
$ws.Range("E3").Value = 15.78412391665828
$ws.Range("F3").Value = 37.86731700481619
$ws.Range("G3").Value = 3.649386480677572
$ws.Range("J3").Value = 7.799749773948999
$ws.Range("L3").Value = 12.64155157016876
$ws.Range("N3").Value = 18.00619248413137
$ws.Range("O3").Value = 24.27136305561313
# Row 4
$ws.Range("B4").Value = 17.04254286890179
$ws.Range("C4").Value = 10.42415064130647
$ws.Range("E4").Value = 15.80150863720164
$ws.Range("F4").Value = 37.87882779522931
$ws.Range("G4").Value = 3.650803676849022
$ws.Range("J4").Value = 7.803818088074102
$ws.Range("L4").Value = 12.62698856031388
$ws.Range("N4").Value = 18.04143792110278
$ws.Range("O4").Value = 24.31208153850937
# Row 5
$ws.Range("B5").Value = 16.9299206556445
$ws.Range("C5").Value = 10.39208234087552
$ws.Range("E5").Value = 15.80910344715583
$ws.Range("F5").Value = 37.88591036140324
$ws.Range("G5").Value = 3.651399279116468
$ws.Range("J5").Value = 7.805545589121916
$ws.Range("L5").Value = 12.62148119094821
$ws.Range("N5").Value = 18.05627265266385
$ws.Range("O5").Value = 24.33019600940542
# Row 6
$ws.Range("B6").Value = 16.91118838502911
$ws.Range("C6").Value = 10.38673774114418
$ws.Range("E6").Value = 15.81039539662991
$ws.Range("F6").Value = 37.88723083444418
$ws.Range("G6").Value = 3.651499272310536
$ws.Range("J6").Value = 7.805836651178873
$ws.Range("L6").Value = 12.62059260876897
$ws.Range("N6").Value = 18.05876448251769
$ws.Range("O6").Value = 24.33329564725429
# Row 7
$ws.Range("B7").Value = 17.04102622382521
$ws.Range("C7").Value = 10.42371948890629
$ws.Range("E7").Value = 15.80160899640937
$ws.Range("F7").Value = 37.87891363050866
$ws.Range("G7").Value = 3.650811636053948
$ws.Range("J7").Value = 7.803841103545306
$ws.Range("L7").Value = 12.62691255128181
$ws.Range("N7").Value = 18.04163607528642
$ws.Range("O7").Value = 24.31231968309989
# Row 8
$ws.Range("B8").Value = 17.60823252497582
$ws.Range("C8").Value = 10.58391706557683
$ws.Range("E8").Value = 15.76731636581758
$ws.Range("F8").Value = 37.86321148163442
$ws.Range("G8").Value = 3.647935849324307
$ws.Range("J8").Value = 7.795645903766555
$ws.Range("L8").Value = 12.65854059330554
$ws.Range("N8").Value = 17.97019255768919
$ws.Range("O8").Value = 24.23312610454913
# Row 9
$ws.Range("B9").Value = 18.68942172066875
$ws.Range("C9").Value = 10.88606030266115
$ws.Range("E9").Value = 15.71614388970119
$ws.Range("F9").Value = 37.90786806255323
$ws.Range("G9").Value = 3.642862576187543
$ws.Range("J9").Value = 7.781755807451822
$ws.Range("L9").Value = 12.73372745401245
$ws.Range("N9").Value = 17.84490429130404
$ws.Range("O9").Value = 24.12616995717184
# Row 10
$ws.Range("B10").Value = 19.45627924736428
$ws.Range("C10").Value = 11.09966240315265
$ws.Range("E10").Value = 15.68832239918837
$ws.Range("F10").Value = 37.98653356578144
$ws.Range("G10").Value = 3.639476677735979
$ws.Range("J10").Value = 7.772866988131337
$ws.Range("L10").Value = 12.79672537505263
$ws.Range("N10").Value = 17.76182030418322
$ws.Range("O10").Value = 24.07729806472737
# Row 11
$ws.Range("B11").Value = 19.79756325057667
$ws.Range("C11").Value = 11.1947994014201
$ws.Range("E11").Value = 15.67778202751107
$ws.Range("F11").Value = 38.0322162524343
$ws.Range("G11").Value = 3.638009697974228
$ws.Range("J11").Value = 7.769106248209062
$ws.Range("L11").Value = 12.8270068196015
$ws.Range("N11").Value = 17.72595656238238
$ws.Range("O11").Value = 24.06155784487707
# Row 12
$ws.Range("B12").Value = 19.92559052978958
$ws.Range("C12").Value = 11.23051570191954
$ws.Range("E12").Value = 15.67409429490798
$ws.Range("F12").Value = 38.05093041593203
$ws.Range("G12").Value = 3.637464669163929
$ws.Range("J12").Value = 7.767722602746203
$ws.Range("L12").Value = 12.8387013006399
$ws.Range("N12").Value = 17.7126526903986
$ws.Range("O12").Value = 24.05653356318529
# Row 13
$ws.Range("B13").Value = 19.89807319477484
$ws.Range("C13").Value = 11.22283762337032
$ws.Range("E13").Value = 15.67487501758392
$ws.Range("F13").Value = 38.0468371940876
$ws.Range("G13").Value = 3.637581585410211
$ws.Range("J13").Value = 7.768018799117557
$ws.Range("L13").Value = 12.83617266068324
$ws.Range("N13").Value = 17.71550561253955
$ws.Range("O13").Value = 24.05757395731013
# Row 14
$ws.Range("B14").Value = 19.80812103604119
$ws.Range("C14").Value = 11.19774411623505
$ws.Range("E14").Value = 15.6774725528239
$ws.Range("F14").Value = 38.03372756888027
$ws.Range("G14").Value = 3.637964648303243
$ws.Range("J14").Value = 7.768991605057018
$ws.Range("L14").Value = 12.82796440564376
$ws.Range("N14").Value = 17.72485649891791
$ws.Range("O14").Value = 24.06112571788475
# Row 15
$ws.Range("B15").Value = 19.75286177331601
$ws.Range("C15").Value = 11.18233272581036
$ws.Range("E15").Value = 15.67910314816566
$ws.Range("F15").Value = 38.02588156771399
$ws.Range("G15").Value = 3.638200649212525
$ws.Range("J15").Value = 7.76959274045474
$ws.Range("L15").Value = 12.82296606790552
$ws.Range("N15").Value = 17.73062023008923
$ws.Range("O15").Value = 24.0634232651224
# Row 16
$ws.Range("B16").Value = 19.43381268374324
$ws.Range("C16").Value = 11.09340255292394
$ws.Range("E16").Value = 15.6890537646322
$ws.Range("F16").Value = 37.98374659837946
$ws.Range("G16").Value = 3.639574018374665
$ws.Range("J16").Value = 7.773118435002425
$ws.Range("L16").Value = 12.79477858068067
$ws.Range("N16").Value = 17.76420288257303
$ws.Range("O16").Value = 24.07845759829265
# Row 17
$ws.Range("B17").Value = 19.23605924500937
$ws.Range("C17").Value = 11.0383135068671
$ws.Range("E17").Value = 15.69569965734562
$ws.Range("F17").Value = 37.96042784934798
$ws.Range("G17").Value = 3.640435267179698
$ws.Range("J17").Value = 7.775353627353711
$ws.Range("L17").Value = 12.77789814987148
$ws.Range("N17").Value = 17.78529889202881
$ws.Range("O17").Value = 24.08934546001281
# Row 18
$ws.Range("B18").Value = 19.12161144442928
$ws.Range("C18").Value = 11.00643754465631
$ws.Range("E18").Value = 15.69972138210966
$ws.Range("F18").Value = 37.94794769426733
$ws.Range("G18").Value = 3.640937535563718
$ws.Range("J18").Value = 7.776665884260741
$ws.Range("L18").Value = 12.76834214080488
$ws.Range("N18").Value = 17.79761461506627
$ws.Range("O18").Value = 24.09621878528657
# Row 19
$ws.Range("B19").Value = 19.08274414050686
$ws.Range("C19").Value = 10.9956127789808
$ws.Range("E19").Value = 15.70111729527128
$ws.Range("F19").Value = 37.94388245060023
$ws.Range("G19").Value = 3.641108781868031
$ws.Range("J19").Value = 7.777114772480346
$ws.Range("L19").Value = 12.7651331222418
$ws.Range("N19").Value = 17.80181576638447
$ws.Range("O19").Value = 24.09865081049397
# Row 20
$ws.Range("B20").Value = 19.25718434285037
$ws.Range("C20").Value = 11.04419764623255
$ws.Range("E20").Value = 15.69497157996022
$ws.Range("F20").Value = 37.96281374724811
$ws.Range("G20").Value = 3.640342871949338
$ws.Range("J20").Value = 7.775112932341574
$ws.Range("L20").Value = 12.77967929003011
$ws.Range("N20").Value = 17.78303437029129
$ws.Range("O20").Value = 24.08812318006406
# Row 21
$ws.Range("B21").Value = 19.83457591465451
$ws.Range("C21").Value = 11.20512324056164
$ws.Range("E21").Value = 15.67670135683105
$ws.Range("F21").Value = 38.03753985009389
$ws.Range("G21").Value = 3.637851849319888
$ws.Range("J21").Value = 7.768704771792392
$ws.Range("L21").Value = 12.83036924147367
$ws.Range("N21").Value = 17.72210240694739
$ws.Range("O21").Value = 24.06005705276655
# Row 22
$ws.Range("B22").Value = 20.20483492834601
$ws.Range("C22").Value = 11.3084825711347
$ws.Range("E22").Value = 15.66653047226222
$ws.Range("F22").Value = 38.09462120762724
$ws.Range("G22").Value = 3.636284909775687
$ws.Range("J22").Value = 7.764752435306882
$ws.Range("L22").Value = 12.86482151734685
$ws.Range("N22").Value = 17.68389375681788
$ws.Range("O22").Value = 24.04717197098064
# Row 23
$ws.Range("B23").Value = 20.00790804220272
$ws.Range("C23").Value = 11.2534895481849
$ws.Range("E23").Value = 15.67179713285371
$ws.Range("F23").Value = 38.06340461827155
$ws.Range("G23").Value = 3.637115642883113
$ws.Range("J23").Value = 7.766840366667993
$ws.Range("L23").Value = 12.84631460041886
$ws.Range("N23").Value = 17.70413902097063
$ws.Range("O23").Value = 24.05354884540179
# Row 24
$ws.Range("B24").Value = 19.24763603370535
$ws.Range("C24").Value = 11.04153806257904
$ws.Range("E24").Value = 15.69530011788955
$ws.Range("F24").Value = 37.96173219743613
$ws.Range("G24").Value = 3.640384621634718
$ws.Range("J24").Value = 7.775221665761597
$ws.Range("L24").Value = 12.77887357270687
$ws.Range("N24").Value = 17.78405757691541
$ws.Range("O24").Value = 24.08867386120849
# Row 25
$ws.Range("B25").Value = 18.40117938608002
$ws.Range("C25").Value = 10.80571938863249
$ws.Range("E25").Value = 15.72826881998118
$ws.Range("F25").Value = 37.88771820128618
$ws.Range("G25").Value = 3.644174805860693
$ws.Range("J25").Value = 7.785281353593398
$ws.Range("L25").Value = 12.71200453864012
$ws.Range("N25").Value = 17.84490429130404
$ws.Range("O25").Value = 24.12616995717184
